# Update loading_percent values for case "380 kV" (Case_2_110)
# Updates columns B, C, E, F, G, I, J, N for rows 2-25 (data index 0-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data block: row|B,C,E,F,G,I,J,N (new values, in column order)
$dataText = @"
2|17.69814383710417,14.17583283824367,18.861019024475,45.16944210495318,3.652468423574997,23.679235715812,8.765655687048856,17.3715774693497
3|17.09024421006961,13.54621278142966,18.74643778654163,44.81053405387173,3.656613830622018,23.70178474691726,8.790270211142079,17.4529703844758
4|16.71100213583595,13.14889310772383,18.67934522031161,44.60358811666204,3.659287199767737,23.72339157788267,8.807083565012116,17.5049977256596
5|16.55521771192892,12.98454776540675,18.65284510658708,44.52270234611378,3.660408958357569,23.73413485232052,8.814361249709153,17.52671743707786
6|16.52928310200349,12.95711968780527,18.64849613614708,44.50948139377181,3.660597182641798,23.73603541987035,8.81559539952576,17.53035533828027
7|16.70890582749505,13.14668617214669,18.6789844003344,44.60248322125759,3.659302197082657,23.72352863609142,8.807179990944803,17.50528854448914
8|17.48993603423393,13.96111507228201,18.820847435015,45.04294878288931,3.65387125943723,23.68539197225244,8.773789141755742,17.39921724556921
9|18.96342483530048,15.46310713903292,19.12394834303801,46.00978039998996,3.644231276189285,23.67273783947384,8.721858524462698,17.20738840458037
10|19.99794008310855,16.49739815115373,19.36045671275163,46.77790242686034,3.637755950377945,23.70198556551735,8.692047355972445,17.07616607004968
11|20.455893811054,16.95104865276532,19.47074146310058,47.13870324287387,3.634940154893422,23.72377583997388,8.680314326858459,17.018547181786
12|20.6273304565016,17.12028735181184,19.51286381781616,47.27685944695918,3.633892418000416,23.73325451597809,8.67613562629168,16.99702435223323
13|20.59049895433366,17.08395392120804,19.5037764556672,47.24703878629811,3.634117244221764,23.73115843911478,8.677023804718074,17.00164653981132
14|20.47003855439972,16.96502375442014,19.47419984234162,47.150039424211,3.634853586051119,23.72453102724108,8.679965236620404,17.01677056188836
15|20.39599070614943,16.89184026729993,19.45612935061667,47.09082036433993,3.635307027758969,23.72063154227421,8.681801416375258,17.02607297098693
16|19.96774292281074,16.46740071135632,19.35330123432584,46.7545430952219,3.637942571257622,23.70073289783923,8.692851049733017,17.07997316554065
17|19.70166544434649,16.20260852507992,19.29089044600544,46.5510867437842,3.639592559976986,23.69070446419507,8.700098930149457,17.11356905899605
18|19.54744217479377,16.04872684659643,19.25524886573593,46.43514166787413,3.640553820926148,23.68573472950757,8.704439715059259,17.13308794856875
19|19.49502694638274,15.99635790073852,19.24322594951945,46.39607296206493,3.640881392038216,23.68418894365913,8.705938926190383,17.13973033864064
20|19.73011346224309,16.23096062303677,19.29750794160763,46.57263425658643,3.639415650936594,23.69168932604449,8.699309569442661,17.10997250920589
21|20.50547563368335,17.00002650548894,19.48287766665006,47.17848986891084,3.634636802372126,23.72644430529056,8.679094081166214,17.01232024862691
22|21.0006141841642,17.48775379422296,19.6061121080851,47.58331057301563,3.631621575329082,23.75631378132484,8.66742353773709,16.95022428690577
23|20.73745970238961,17.22884451764818,19.54015813263716,47.36647576833267,3.633221017321091,23.73971526340601,8.673510812147141,16.98320890715643
24|19.71725599995122,16.21814776648477,19.29451542465223,46.56288942625642,3.639495592069371,23.69124159198062,8.699665897924568,17.11159787328068
25|18.57241684777333,15.06822823281613,19.039421827578,45.73769007455616,3.64673190157322,23.66943536515696,8.734448045679033,17.25756660788098
"@

$columns = @("B", "C", "E", "F", "G", "I", "J", "N")

$lines = $dataText -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split "\|"
    $row = [int]$parts[0]
    $values = $parts[1] -split ","
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $val = [double]$values[$i]
        $ws.Range("$col$row").Value = $val
    }
}
